# Update daily stats sheet: ut 09. 02. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing AgTests (H) / AgPosit (I) values for prior days ---
$ws.Range("H317").Value = 61383

$ws.Range("H320").Value = 76615

$ws.Range("H321").Value = 90733

$ws.Range("H322").Value = 107038

$ws.Range("H323").Value = 148776
$ws.Range("I323").Value = 2285

$ws.Range("H324").Value = 232565
$ws.Range("I324").Value = 2659

$ws.Range("H325").Value = 725776

$ws.Range("H326").Value = 425491

$ws.Range("H327").Value = 236369
$ws.Range("I327").Value = 2875

$ws.Range("H328").Value = 178247
$ws.Range("I328").Value = 2617

$ws.Range("H329").Value = 82243

$ws.Range("H330").Value = 70766

$ws.Range("H331").Value = 149227
$ws.Range("I331").Value = 2580

$ws.Range("H332").Value = 417199
$ws.Range("I332").Value = 4073

$ws.Range("H333").Value = 256303
$ws.Range("I333").Value = 2737

$ws.Range("H334").Value = 207048
$ws.Range("I334").Value = 3442

$ws.Range("H335").Value = 129541
$ws.Range("I335").Value = 2898

$ws.Range("H336").Value = 100300
$ws.Range("I336").Value = 3178

$ws.Range("H337").Value = 102274
$ws.Range("I337").Value = 2882

$ws.Range("H338").Value = 214913
$ws.Range("I338").Value = 3763

$ws.Range("H339").Value = 563204
$ws.Range("I339").Value = 4368

$ws.Range("H340").Value = 310254
$ws.Range("I340").Value = 2705

# --- Append new row 341 with the latest day's data ---
$ws.Range("A341").Value = 44235
$ws.Range("A341").NumberFormat = "yyyy-mm-dd"
$ws.Range("B341").Value = 265807
$ws.Range("C341").Value = 248386
$ws.Range("D341").Value = 12039
$ws.Range("E341").Value = 9410
$ws.Range("F341").Value = 1724
$ws.Range("G341").Value = 5382
$ws.Range("H341").Value = 370937
$ws.Range("I341").Value = 4474
